$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-Rows $ws 135 136
Swap-Rows $ws 152 153
